$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wait Times")

# Step 1: set M column values in the precise order needed so new shared strings are appended in target order
$ws.Range("M78").Value = 'MBL lab not returning data on time; samples backlogged'
$ws.Range("M79").Value = 'MBL lab not returning data on time; samples backlogged'
$ws.Range("M80").Value = 'MBL lab not returning data on time; samples backlogged'
$ws.Range("M81").Value = 'MBL lab not returning data on time; samples backlogged'
$ws.Range("M106").Value = 'MBL lab not returning data on time; samples backlogged'
$ws.Range("M107").Value = 'MBL lab not returning data on time; samples backlogged'
$ws.Range("M108").Value = 'MBL lab not returning data on time; samples backlogged'
$ws.Range("M109").Value = 'MBL lab not returning data on time; samples backlogged'
$ws.Range("M110").Value = 'MBL lab not returning data on time; samples backlogged'
$ws.Range("M61").Value = 'MBL lab not returning soil data on time; samples backlogged; dependency'
$ws.Range("M62").Value = 'MBL lab not returning soil data on time; samples backlogged; dependency'
$ws.Range("M63").Value = 'MBL lab not returning soil data on time; samples backlogged; dependency'
$ws.Range("M76").Value = 'NRCS not returning data on time; samples backlogged'
$ws.Range("M97").Value = 'NRCS not returning data on time'
$ws.Range("M98").Value = 'NRCS not returning data on time'
$ws.Range("M99").Value = 'NRCS not returning data on time'
$ws.Range("M100").Value = 'NRCS not returning data on time'
$ws.Range("M101").Value = 'Waiting on HQ Scientist availability to complete error correction of 2013 data'
$ws.Range("M102").Value = 'Waiting on HQ Scientist availability to complete error correction of 2013 data'
$ws.Range("M103").Value = 'Waiting on HQ Scientist availability to complete error correction of 2013 data'
$ws.Range("M104").Value = 'Waiting on HQ Scientist availability to complete error correction of 2013 data'
$ws.Range("M105").Value = 'Waiting on HQ Scientist availability to complete error correction of 2013 data'
$ws.Range("M112").Value = 'Awaiting 2017 data - proceeding on schedule; requires final sampling bout for year to be completed; all ticks identified to species; then samples selected for pathogen testing'
$ws.Range("M113").Value = 'Awaiting 2017 data - proceeding on schedule; requires final sampling bout for year to be completed; all ticks identified to species; then samples selected for pathogen testing'

# Step 2: set remaining M column values (reusing existing shared strings)
$ws.Range("M68").Value = 'U WY lab not returning data on time; samples backlogged'
$ws.Range("M72").Value = 'Vendor identified; contracting process prolonged'
$ws.Range("M73").Value = 'Vendor identified; contracting process prolonged'
$ws.Range("M74").Value = 'Vendor identified; contracting process prolonged'
$ws.Range("M75").Value = 'Vendor identified; contracting process prolonged'
$ws.Range("M77").Value = 'U WY lab not returning data on time; samples backlogged'
$ws.Range("M111").Value = 'U WY lab not returning data on time; samples backlogged'
$ws.Range("M114").Value = 'Waiting on HQ Scientist availability to conduct QA and upload 2017 data'
$ws.Range("M115").Value = 'Waiting on HQ Scientist availability to conduct QA and upload 2017 data'
$ws.Range("M116").Value = 'Waiting on HQ Scientist availability to conduct QA and upload 2017 data'

# Step 3: set J (numeric) and K (text) columns for all rows
$ws.Range("J61").Value = 80
$ws.Range("K61").Value = 'external lab processing'
$ws.Range("J62").Value = 80
$ws.Range("K62").Value = 'external lab processing'
$ws.Range("J63").Value = 80
$ws.Range("K63").Value = 'external lab processing'
$ws.Range("J64").Value = 100
$ws.Range("K64").Value = 'NA'
$ws.Range("J65").Value = 100
$ws.Range("K65").Value = 'NA'
$ws.Range("J66").Value = 100
$ws.Range("K66").Value = 'NA'
$ws.Range("J67").Value = 100
$ws.Range("K67").Value = 'NA'
$ws.Range("J68").Value = 50
$ws.Range("K68").Value = 'external lab processing'
$ws.Range("J72").Value = 6
$ws.Range("K72").Value = 'contracting'
$ws.Range("J73").Value = 6
$ws.Range("K73").Value = 'contracting'
$ws.Range("J74").Value = 6
$ws.Range("K74").Value = 'contracting'
$ws.Range("J75").Value = 6
$ws.Range("K75").Value = 'contracting'
$ws.Range("J76").Value = 15
$ws.Range("K76").Value = 'external lab processing'
$ws.Range("J77").Value = 95
$ws.Range("K77").Value = 'external lab processing'
$ws.Range("J78").Value = 75
$ws.Range("K78").Value = 'external lab processing'
$ws.Range("J79").Value = 75
$ws.Range("K79").Value = 'external lab processing'
$ws.Range("J80").Value = 75
$ws.Range("K80").Value = 'external lab processing'
$ws.Range("J81").Value = 75
$ws.Range("K81").Value = 'external lab processing'
$ws.Range("J82").Value = 70
$ws.Range("K82").Value = 'external lab processing'
$ws.Range("J83").Value = 70
$ws.Range("K83").Value = 'external lab processing'
$ws.Range("J84").Value = 70
$ws.Range("K84").Value = 'external lab processing'
$ws.Range("J85").Value = 1
$ws.Range("K85").Value = 'external lab processing'
$ws.Range("J86").Value = 1
$ws.Range("K86").Value = 'external lab processing'
$ws.Range("J87").Value = 2
$ws.Range("K87").Value = 'external lab processing'
$ws.Range("J88").Value = 5
$ws.Range("K88").Value = 'external lab processing'
$ws.Range("J89").Value = 5
$ws.Range("K89").Value = 'external lab processing'
$ws.Range("J90").Value = 6
$ws.Range("K90").Value = 'external lab processing'
$ws.Range("J91").Value = 6
$ws.Range("K91").Value = 'external lab processing'
$ws.Range("J92").Value = 6
$ws.Range("K92").Value = 'external lab processing'
$ws.Range("J93").Value = 6
$ws.Range("K93").Value = 'external lab processing'
$ws.Range("J94").Value = 6
$ws.Range("K94").Value = 'external lab processing'
$ws.Range("J95").Value = 5
$ws.Range("K95").Value = 'external lab processing'
$ws.Range("J96").Value = 5
$ws.Range("K96").Value = 'external lab processing'
$ws.Range("J97").Value = 15
$ws.Range("K97").Value = 'external lab processing'
$ws.Range("J98").Value = 15
$ws.Range("K98").Value = 'external lab processing'
$ws.Range("J99").Value = 15
$ws.Range("K99").Value = 'external lab processing'
$ws.Range("J100").Value = 15
$ws.Range("K100").Value = 'external lab processing'
$ws.Range("J101").Value = 95
$ws.Range("K101").Value = 'science resources'
$ws.Range("J102").Value = 95
$ws.Range("K102").Value = 'science resources'
$ws.Range("J103").Value = 95
$ws.Range("K103").Value = 'science resources'
$ws.Range("J104").Value = 95
$ws.Range("K104").Value = 'science resources'
$ws.Range("J105").Value = 95
$ws.Range("K105").Value = 'science resources'
$ws.Range("J106").Value = 75
$ws.Range("K106").Value = 'external lab processing'
$ws.Range("J107").Value = 75
$ws.Range("K107").Value = 'external lab processing'
$ws.Range("J108").Value = 75
$ws.Range("K108").Value = 'external lab processing'
$ws.Range("J109").Value = 75
$ws.Range("K109").Value = 'external lab processing'
$ws.Range("J110").Value = 75
$ws.Range("K110").Value = 'external lab processing'
$ws.Range("J111").Value = 95
$ws.Range("K111").Value = 'external lab processing'
$ws.Range("J112").Value = 70
$ws.Range("K112").Value = 'external lab processing'
$ws.Range("J113").Value = 75
$ws.Range("K113").Value = 'external lab processing'
$ws.Range("J114").Value = 75
$ws.Range("K114").Value = 'science resources'
$ws.Range("J115").Value = 75
$ws.Range("K115").Value = 'science resources'
$ws.Range("J116").Value = 75
$ws.Range("K116").Value = 'science resources'

# Step 4: set L column - either literal NA text, date (copy style from existing date cell then set value), or blank yellow style
$ws.Range("L64").Value = 'NA'
$ws.Range("L65").Value = 'NA'
$ws.Range("L66").Value = 'NA'
$ws.Range("L67").Value = 'NA'

# Step 4b: date cells - copy number format style from L12 (existing date-styled cell) then set serial value
$ws.Range("L12").Copy()
$ws.Range("L61").PasteSpecial(-4122)
$ws.Range("L61").Value = 43344
$ws.Range("L62").PasteSpecial(-4122)
$ws.Range("L62").Value = 43344
$ws.Range("L63").PasteSpecial(-4122)
$ws.Range("L63").Value = 43344
$ws.Range("L68").PasteSpecial(-4122)
$ws.Range("L68").Value = 43344
$ws.Range("L72").PasteSpecial(-4122)
$ws.Range("L72").Value = 43374
$ws.Range("L73").PasteSpecial(-4122)
$ws.Range("L73").Value = 43374
$ws.Range("L74").PasteSpecial(-4122)
$ws.Range("L74").Value = 43374
$ws.Range("L75").PasteSpecial(-4122)
$ws.Range("L75").Value = 43374
$ws.Range("L76").PasteSpecial(-4122)
$ws.Range("L76").Value = 43465
$ws.Range("L77").PasteSpecial(-4122)
$ws.Range("L77").Value = 43344
$ws.Range("L78").PasteSpecial(-4122)
$ws.Range("L78").Value = 43344
$ws.Range("L79").PasteSpecial(-4122)
$ws.Range("L79").Value = 43344
$ws.Range("L80").PasteSpecial(-4122)
$ws.Range("L80").Value = 43344
$ws.Range("L81").PasteSpecial(-4122)
$ws.Range("L81").Value = 43344
$ws.Range("L97").PasteSpecial(-4122)
$ws.Range("L97").Value = 43465
$ws.Range("L98").PasteSpecial(-4122)
$ws.Range("L98").Value = 43465
$ws.Range("L99").PasteSpecial(-4122)
$ws.Range("L99").Value = 43465
$ws.Range("L100").PasteSpecial(-4122)
$ws.Range("L100").Value = 43465
$ws.Range("L101").PasteSpecial(-4122)
$ws.Range("L101").Value = 43251
$ws.Range("L102").PasteSpecial(-4122)
$ws.Range("L102").Value = 43251
$ws.Range("L103").PasteSpecial(-4122)
$ws.Range("L103").Value = 43251
$ws.Range("L104").PasteSpecial(-4122)
$ws.Range("L104").Value = 43251
$ws.Range("L105").PasteSpecial(-4122)
$ws.Range("L105").Value = 43251
$ws.Range("L106").PasteSpecial(-4122)
$ws.Range("L106").Value = 43344
$ws.Range("L107").PasteSpecial(-4122)
$ws.Range("L107").Value = 43344
$ws.Range("L108").PasteSpecial(-4122)
$ws.Range("L108").Value = 43344
$ws.Range("L109").PasteSpecial(-4122)
$ws.Range("L109").Value = 43344
$ws.Range("L110").PasteSpecial(-4122)
$ws.Range("L110").Value = 43344
$ws.Range("L111").PasteSpecial(-4122)
$ws.Range("L111").Value = 43344
$ws.Range("L112").PasteSpecial(-4122)
$ws.Range("L112").Value = 43344
$ws.Range("L113").PasteSpecial(-4122)
$ws.Range("L113").Value = 43344
$ws.Range("L114").PasteSpecial(-4122)
$ws.Range("L114").Value = 43343
$ws.Range("L115").PasteSpecial(-4122)
$ws.Range("L115").Value = 43343
$ws.Range("L116").PasteSpecial(-4122)
$ws.Range("L116").Value = 43343

# Step 5: blank yellow-highlight style cells for L and M in rows with no date/note (82-96 excluded 85-96 set individually)
$ws.Range("L82:M82").Interior.Color = 65535
$ws.Range("L83:M83").Interior.Color = 65535
$ws.Range("L84:M84").Interior.Color = 65535
$ws.Range("L85:M85").Interior.Color = 65535
$ws.Range("L86:M86").Interior.Color = 65535
$ws.Range("L87:M87").Interior.Color = 65535
$ws.Range("L88:M88").Interior.Color = 65535
$ws.Range("L89:M89").Interior.Color = 65535
$ws.Range("L90:M90").Interior.Color = 65535
$ws.Range("L91:M91").Interior.Color = 65535
$ws.Range("L92:M92").Interior.Color = 65535
$ws.Range("L93:M93").Interior.Color = 65535
$ws.Range("L94:M94").Interior.Color = 65535
$ws.Range("L95:M95").Interior.Color = 65535
$ws.Range("L96:M96").Interior.Color = 65535

# Step 6: restore active selection to match the latest position the author was working at
$ws.Activate()
$ws.Range("O18").Select()
